# Added Configurable zero_before_threshold parameter to enable setting dims
# before noise_threshold or First Rise Point to 0.
#
# This updates the First_Noticeable_Increase_Index (C),
# First_Noticeable_Increase_Cumulative_Value (E), and Pulse_Width (G)
# columns on each of the Step3_DataPts_* sheets to reflect the new
# zero_before_threshold behavior.

$wb = $excel.ActiveWorkbook

$sheetNames = @("Step3_DataPts_0.5", "Step3_DataPts_0.7", "Step3_DataPts_0.8", "Step3_DataPts_0.9")

# New values for column C (First_Noticeable_Increase_Index), keyed by row (2-6).
$newC = @{ 2 = 89; 3 = 87; 4 = 87; 5 = 88; 6 = 89 }

# New values for column E (First_Noticeable_Increase_Cumulative_Value), keyed by row (2-6).
$newE = @{
    2 = 0.01346169945484399
    3 = 0.007240865138925249
    4 = 0.005269976108616141
    5 = 0.01034909697080863
    6 = 0.01054872919713661
}

# New values for column G (Pulse_Width), keyed by sheet name then row (2-6).
$newG = @{
    "Step3_DataPts_0.5" = @{ 2 = 16; 3 = 17; 4 = 17; 5 = 17; 6 = 16 }
    "Step3_DataPts_0.7" = @{ 2 = 29; 3 = 28; 4 = 29; 5 = 28; 6 = 27 }
    "Step3_DataPts_0.8" = @{ 2 = 63; 3 = 72; 4 = 65; 5 = 72; 6 = 71 }
    "Step3_DataPts_0.9" = @{ 2 = 75; 3 = 83; 4 = 80; 5 = 81; 6 = 80 }
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in 2..6) {
        $ws.Cells.Item($row, 3).Value = $newC[$row]
        $ws.Cells.Item($row, 5).Value = $newE[$row]
        $ws.Cells.Item($row, 7).Value = $newG[$sheetName][$row]
    }
}
